$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column C with the "like a pro" article follower counts,
# mirroring the date/format of column B.
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("C1").Value = 44186

$ws.Range("C2").Value = 328
$ws.Range("C3").Value = 458
$ws.Range("C4").Value = 941
$ws.Range("C5").Value = 82

# Move the active selection as in the saved workbook.
$ws.Range("F11").Select()
